$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain text formatting for numeric-looking strings
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "70.119.99"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "3.544.99"
$ws.Range("E3").Value = "  +0.32%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "603.49"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").Value = "197.32"
$ws.Range("E6").Value = "  +6.33%  "
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("E9").Value = "  -2.51%  "
$ws.Range("D10").Value = "0.655"
$ws.Range("E10").Value = "  -0.10%  "
$ws.Range("D11").Value = "54.08"
$ws.Range("E11").Value = "  +1.23%  "
$ws.Range("E12").Value = "  -0.70%  "
$ws.Range("D13").Value = "9.56"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "4.112.34"
$ws.Range("D15").Value = "604.56"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").Value = "19.27"
$ws.Range("E16").Value = "  +1.72%  "
$ws.Range("D17").Value = "70.225.13"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("E18").Value = "  -0.31%  "
$ws.Range("D19").Value = "3.537.27"
$ws.Range("E19").Value = "  -0.50%  "
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").Value = "0.997"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "18.14"
$ws.Range("E22").Value = "  +3.53%  "
$ws.Range("D23").Value = "5.29"
$ws.Range("E23").Value = "  +6.09%  "
$ws.Range("D24").Value = "102.88"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("E26").Value = "  +4.35%  "
$ws.Range("D27").Value = "10.97"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "9.64"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("D29").Value = "33.83"
$ws.Range("E29").Value = "  -0.99%  "
$ws.Range("D30").Value = "4.39"
$ws.Range("E30").Value = "  +21.13%  "
$ws.Range("D31").Value = "7.15"
$ws.Range("E31").Value = "  +1.11%  "
$ws.Range("D32").Value = "12.68"
$ws.Range("E32").Value = "  +2.89%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("D34").Value = "63.43"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").Value = "0.0₃0846"
$ws.Range("E35").Value = "  +8.63%  "
$ws.Range("D36").Value = "3.789.69"
$ws.Range("B37").Value = "Stacks"
$ws.Range("C37").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D37").Value = "3.70"
$ws.Range("E37").Value = "  +3.53%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "3.08"
$ws.Range("E38").Value = "  -4.34%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").Value = "36.84"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("D42").Value = "489.57"
$ws.Range("E42").Value = "  -7.68%  "
$ws.Range("E43").Value = "  -1.40%  "
$ws.Range("D44").Value = "0.0459"
$ws.Range("E44").Value = "  -1.52%  "
$ws.Range("E45").Value = "  -3.24%  "
$ws.Range("E46").Value = "  -2.24%  "
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("D49").Value = "8.69"
$ws.Range("E49").Value = "  -3.95%  "
$ws.Range("D50").Value = "0.000251"
$ws.Range("E50").Value = "  +3.74%  "
$ws.Range("D51").Value = "130.91"
$ws.Range("E51").Value = "  -2.15%  "
